# "added parallel test methods"
# Split the single Sheet1 (with 3 locale hyperlinks stacked in A2:A4) into three
# parallel per-locale sheets: dcPages1 (kr, plain text), dcPages2 (tr, hyperlink),
# dcPages3 (ar, hyperlink).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$krUrl = "https://www.stage.adobe.com/kr/acrobat/pricing/compare-versions.html"
$trUrl = "https://www.stage.adobe.com/tr/acrobat/pricing/compare-versions.html"
$arUrl = "https://www.stage.adobe.com/ar/acrobat/pricing/compare-versions.html"

# ---- dcPages1 (was Sheet1) ------------------------------------------------
$ws1.Name = "dcPages1"

# Drop all the hyperlinks that used to live on A2:A4.
$ws1.Hyperlinks.Delete()

# A2 keeps plain text (the kr url), no hyperlink / no special style anymore.
$ws1.Range("A2").Value = $krUrl
$ws1.Range("A2").Style = "Normal"

# A3/A4 (previously tr/ar) become fully blank again, like the rest of the sheet.
$ws1.Range("A3:A4").Clear()

# Column got a lot wider (no longer bestfit to the short placeholder text) and
# the header row got a touch taller.
$ws1.Columns("A").ColumnWidth = 50.8
$ws1.Rows(1).RowHeight = 19.5

# Selection moves up to A2.
$ws1.Range("A2").Select()

# ---- dcPages2 (new sheet, tr link) ----------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "dcPages2"

$ws2.Range("A1").Value = "url"
$ws2.Range("A1").Style = "Normal"
$ws2.Range("A1").HorizontalAlignment = -4131

$ws2.Range("A2").Value = $trUrl
$ws2.Hyperlinks.Add($ws2.Range("A2"), $trUrl)

$ws2.Range("A2").Select()

# ---- dcPages3 (new sheet, ar link) ----------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "dcPages3"

$ws3.Range("A1").Value = "url"
$ws3.Range("A1").Style = "Normal"
$ws3.Range("A1").HorizontalAlignment = -4131

$ws3.Range("A2").Value = $arUrl
$ws3.Hyperlinks.Add($ws3.Range("A2"), $arUrl)

$ws3.Range("A2").Select()

# Leave dcPages1 as the active/visible tab, matching the original authoring flow.
$ws1.Select()
$ws1.Range("A2").Select()
